$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear existing contents across the full original range
$ws.Range("A1:C25").ClearContents()

# 2. Reset all row heights to default (auto) before re-applying custom ones
$ws.Range("A1:C25").EntireRow.AutoFit()

# 3. Write target cell values
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOM3113"
$ws.Range("C2").Value = "LOM3113"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Tratamentos de Minérios e Hidrometalurgia"
$ws.Range("C3").Value = " Tratamentos de Minérios e Hidrometalurgia"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Tratamento f Ores and Hydrometallurgy"
$ws.Range("C4").Value = "Tratamento f Ores and Hydrometallurgy"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EM-5"
$ws.Range("C9").Value = "EM-5"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"

$ws.Range("A14").Value = "Short syllabus:"

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("A16").Value = "Syllabus:"

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Serão realizadas duas provas escritas (P1 e P2) com peso 1. No mínimo, um relatório a partir de trabalhos em grupo, com peso 1 (NR) e avaliação individual realizada durante todo o curso (AI), com peso 1."
$ws.Range("C19").Value = "Serão realizadas duas provas escritas (P1 e P2) com peso 1. No mínimo, um relatório a partir de trabalhos em grupo, com peso 1 (NR) e avaliação individual realizada durante todo o curso (AI), com peso 1."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final (NF) será calculada pela equação: NF = 0,5[(P1 + P2)/2] + 0,3NR + 0,2AI."
$ws.Range("C20").Value = "A nota final (NF) será calculada pela equação: NF = 0,5[(P1 + P2)/2] + 0,3NR + 0,2AI."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, inclusive com cobrança das competências desenvolvidas, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, inclusive com cobrança das competências desenvolvidas, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."

$ws.Range("A22").Value = "Requisitos:"

$ws.Range("B23").Value = "LOM3037 -  Química Inorgânica  (Requisito)`n"
$ws.Range("C23").Value = "LOM3037 -  Química Inorgânica  (Requisito)`n"

# 4. Apply custom row heights for rows that need them
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

# 5. Remove now-unused trailing rows (24 and 25) so dimension shrinks to C23
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

